$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 12.785714
$ws.Range("I11").Value = 12.785714
$ws.Range("K11").Value = 12.785714
$ws.Range("M11").Value = 127.214286
$ws.Range("H32").Value = 15230.77
$ws.Range("J32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("N32").Value = -12652
$ws.Range("H40").Value = 1789.0555
$ws.Range("I40").Value = 1667.3334
$ws.Range("J40").Value = 1910.7778
$ws.Range("K40").Value = 1667.3334
$ws.Range("L40").Value = 1910.7778
$ws.Range("M40").Value = -1492.3334
$ws.Range("N40").Value = -2260.7778
$ws.Range("H64").Value = 10927.846
$ws.Range("I64").Value = 7945.4287
$ws.Range("K64").Value = 7945.4287
$ws.Range("M64").Value = -7697.4287
$ws.Range("H67").Value = 10927.846
$ws.Range("I67").Value = 7945.4287
$ws.Range("K67").Value = 7945.4287
$ws.Range("M67").Value = -7087.4287
$ws.Range("H98").Value = 553.7727
$ws.Range("J98").Value = 1893
$ws.Range("L98").Value = 1893
$ws.Range("N98").Value = -4889
$ws.Range("H122").Value = 553.7727
$ws.Range("J122").Value = 1893
$ws.Range("L122").Value = 5679
$ws.Range("N122").Value = -10579
$ws.Range("H131").Value = 26050.572
$ws.Range("I131").Value = 3449.8
$ws.Range("J131").Value = 82552.5
$ws.Range("K131").Value = 10349.4
$ws.Range("L131").Value = 247657.5
$ws.Range("M131").Value = -5309.400000000001
$ws.Range("N131").Value = -257737.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 16745.143
$ws.Range("I22").Value = 3054
$ws.Range("J22").Value = 35000
$ws.Range("K22").Value = 3054
$ws.Range("L22").Value = 35000
$ws.Range("M22").Value = -2755
$ws.Range("N22").Value = -35598
$ws.Range("H97").Value = 1809.05
$ws.Range("I97").Value = 1953.8
$ws.Range("J97").Value = 1374.8
$ws.Range("K97").Value = 1953.8
$ws.Range("L97").Value = 1374.8
$ws.Range("M97").Value = -1457.8
$ws.Range("N97").Value = -2366.8
$ws.Range("H122").Value = 1651.4
$ws.Range("I122").Value = 1106.0834
$ws.Range("K122").Value = 3318.2502
$ws.Range("M122").Value = -868.2501999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1200.7317
$ws.Range("I94").Value = 1068.303
$ws.Range("K94").Value = 1068.303
$ws.Range("M94").Value = -617.3030000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1381.3
$ws.Range("I122").Value = 1002.75
$ws.Range("K122").Value = 3008.25
$ws.Range("M122").Value = -558.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.272728
$ws.Range("I2").Value = 61.25
$ws.Range("K2").Value = 367.5
$ws.Range("M2").Value = -254.5
$ws.Range("H80").Value = 3642
$ws.Range("J80").Value = 5403
$ws.Range("L80").Value = 16209
$ws.Range("N80").Value = -18081
$ws.Range("H83").Value = 3642
$ws.Range("J83").Value = 5403
$ws.Range("L83").Value = 48627
$ws.Range("N83").Value = -57987
$ws.Range("H122").Value = 3267
$ws.Range("J122").Value = 4750.5
$ws.Range("L122").Value = 42754.5
$ws.Range("N122").Value = -47654.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 987.625
$ws.Range("I13").Value = 739
$ws.Range("J13").Value = 1402
$ws.Range("K13").Value = 739
$ws.Range("L13").Value = 1402
$ws.Range("M13").Value = -600
$ws.Range("N13").Value = -1680
$ws.Range("H17").Value = 269.42856
$ws.Range("I17").Value = 519
$ws.Range("J17").Value = 169.6
$ws.Range("K17").Value = 519
$ws.Range("L17").Value = 169.6
$ws.Range("M17").Value = -351
$ws.Range("N17").Value = -505.6
$ws.Range("H102").Value = 2442.5881
$ws.Range("I102").Value = 2470.25
$ws.Range("K102").Value = 2470.25
$ws.Range("M102").Value = -848.25
$ws.Range("H126").Value = 10941.975
$ws.Range("I126").Value = 15075.84
$ws.Range("K126").Value = 45227.52
$ws.Range("M126").Value = -42757.52

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13699.9
$ws.Range("I7").Value = 14875
$ws.Range("J7").Value = 8999.5
$ws.Range("K7").Value = 14875
$ws.Range("L7").Value = 8999.5
$ws.Range("M7").Value = -14763
$ws.Range("N7").Value = -9223.5
$ws.Range("H16").Value = 5471
$ws.Range("I16").Value = 2939
$ws.Range("J16").Value = 6161.5454
$ws.Range("K16").Value = 2939
$ws.Range("L16").Value = 6161.5454
$ws.Range("M16").Value = -2769
$ws.Range("N16").Value = -6501.5454
$ws.Range("H40").Value = 5093.4165
$ws.Range("I40").Value = 4421.222
$ws.Range("K40").Value = 4421.222
$ws.Range("M40").Value = -4285.222
$ws.Range("H46").Value = 2122.111
$ws.Range("I46").Value = 1971.2858
$ws.Range("J46").Value = 2650
$ws.Range("K46").Value = 1971.2858
$ws.Range("L46").Value = 2650
$ws.Range("M46").Value = -1783.2858
$ws.Range("N46").Value = -3026
$ws.Range("H55").Value = 323.85
$ws.Range("I55").Value = 328.41666
$ws.Range("J55").Value = 317
$ws.Range("K55").Value = 328.41666
$ws.Range("L55").Value = 317
$ws.Range("M55").Value = -155.41666
$ws.Range("N55").Value = -663
$ws.Range("H86").Value = 65195
$ws.Range("J86").Value = 65195
$ws.Range("L86").Value = 65195
$ws.Range("N86").Value = -67567
$ws.Range("H89").Value = 65195
$ws.Range("J89").Value = 65195
$ws.Range("L89").Value = 195585
$ws.Range("N89").Value = -207441
$ws.Range("H122").Value = 5044.154
$ws.Range("I122").Value = 4396.5557
$ws.Range("K122").Value = 13189.6671
$ws.Range("M122").Value = -10739.6671
$ws.Range("H126").Value = 13699.9
$ws.Range("I126").Value = 14875
$ws.Range("J126").Value = 8999.5
$ws.Range("K126").Value = 44625
$ws.Range("L126").Value = 26998.5
$ws.Range("M126").Value = -42155
$ws.Range("N126").Value = -31938.5
$ws.Range("H136").Value = 6083.7896
$ws.Range("I136").Value = 5276.5386
$ws.Range("K136").Value = 15829.6158
$ws.Range("M136").Value = -13279.6158

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3654.875
$ws.Range("I122").Value = 3444.4614
$ws.Range("K122").Value = 10333.3842
$ws.Range("M122").Value = -7883.3842
$ws.Range("H124").Value = 30299.666
$ws.Range("J124").Value = 30299.666
$ws.Range("L124").Value = 30299.666
$ws.Range("N124").Value = -40119.666
$ws.Range("H132").Value = 4137.4814
$ws.Range("I132").Value = 3861.652
$ws.Range("K132").Value = 11584.956
$ws.Range("M132").Value = -9054.956
